# Manual de Instalación V4.0.1 - edit script
#
# Change: the instruction paragraph that used to read
#   "Sobre la ventana de CMD abierta inicialmente, correr los siguientes
#    comandos *UNO POR UNO*"
# now reads
#   "Abrir una nueva ventana de CMD COMO ADMINISTRADOR y correr los
#    siguientes comandos *UNO POR UNO*"
# and the whole sentence (including the leading word, which previously was
# not bold) is now bold.

$d = $word.ActiveDocument

$old = "Sobre la ventana de CMD abierta inicialmente, correr los siguientes comandos *UNO POR UNO*"
$new = "Abrir una nueva ventana de CMD COMO ADMINISTRADOR y correr los siguientes comandos *UNO POR UNO*"

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)

# Make sure the whole replaced sentence is bold (previously the leading
# "S" of "Sobre" was not bold; the new leading words "Abrir una nueva"
# must be bold too).
$range2 = $d.Content
$range2.Find.ClearFormatting()
$found2 = $range2.Find.Execute($new, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$range2.Font.Bold = 1
$range2.Font.BoldBi = 1
